$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price cells so Excel does not
# coerce them into floating point numbers (original values are stored as text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '36.357.50'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.936.19'
$ws.Range("E3").Value = '  -2.25%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '241.33'
$ws.Range("E5").Value = '  -1.59%  '
$ws.Range("D6").Value = '0.608'
$ws.Range("E6").Value = '  -3.27%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '56.39'
$ws.Range("E8").Value = '  -5.32%  '
$ws.Range("E9").Value = '  -4.89%  '
$ws.Range("D10").Value = '0.0839'
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("E11").Value = '  -1.45%  '
$ws.Range("D12").Value = '2.220.00'
$ws.Range("E12").Value = '  -2.14%  '
$ws.Range("D13").Value = '0.800'
$ws.Range("E13").Value = '  -7.31%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '20.90'
$ws.Range("E14").Value = '  -11.09%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '13.32'
$ws.Range("E15").Value = '  -4.85%  '
$ws.Range("D16").Value = '5.12'
$ws.Range("E16").Value = '  -6.33%  '
$ws.Range("D17").Value = '1.938.71'
$ws.Range("E17").Value = '  -2.06%  '
$ws.Range("D18").Value = '36.314.98'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").Value = '68.79'
$ws.Range("E19").Value = '  -1.98%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0860'
$ws.Range("E20").Value = '  -1.88%  '
$ws.Range("D21").Value = '226.32'
$ws.Range("E21").Value = '  -3.33%  '
$ws.Range("D22").Value = '4.92'
$ws.Range("E22").Value = '  -7.62%  '
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("E24").Value = '  -10.29%  '
$ws.Range("E25").Value = '  -2.61%  '
$ws.Range("D26").Value = '9.20'
$ws.Range("E26").Value = '  -7.94%  '
$ws.Range("D27").Value = '160.37'
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("D28").Value = '0.131'
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '19.07'
$ws.Range("E29").Value = '  -3.86%  '
$ws.Range("E30").Value = '  -2.74%  '
$ws.Range("D31").Value = '1.10'
$ws.Range("E31").Value = '  -6.88%  '
$ws.Range("E32").Value = '  -7.66%  '
$ws.Range("E33").Value = '  -4.26%  '
$ws.Range("E34").Value = '  -6.39%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").Value = '6.04'
$ws.Range("E36").Value = '  -2.72%  '
$ws.Range("E37").Value = '  +0.34%  '
$ws.Range("D38").Value = '2.12'
$ws.Range("E38").Value = '  -6.38%  '
$ws.Range("D39").Value = '2.95'
$ws.Range("E39").Value = '  -2.21%  '
$ws.Range("D40").Value = '0.0964'
$ws.Range("E40").Value = '  -0.48%  '
$ws.Range("E41").Value = '  -0.85%  '
$ws.Range("E42").Value = '  -3.07%  '
$ws.Range("E43").Value = '  -7.83%  '
$ws.Range("D44").Value = '15.41'
$ws.Range("E44").Value = '  -4.99%  '
$ws.Range("D45").Value = '1.327.39'
$ws.Range("E45").Value = '  -2.86%  '
$ws.Range("E46").Value = '  -7.52%  '
$ws.Range("D47").Value = '85.42'
$ws.Range("E47").Value = '  -7.54%  '
$ws.Range("D48").Value = '7.03'
$ws.Range("E48").Value = '  -6.14%  '
$ws.Range("D49").Value = '2.83'
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").Value = '43.66'
$ws.Range("E50").Value = '  -4.43%  '
$ws.Range("D51").Value = '2.111.18'
$ws.Range("E51").Value = '  -2.18%  '
